$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 already carries cellXf index 1 (the "translated" look used across row 1 / column A)
$styleSource = $ws.Range("A1")

# --- Row 1: M1 ("es") and N1 ("he_IL") pick up the same style, text unchanged ---
$styleSource.Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)

# --- Rows 2-4: he_IL translation -> every touched cell's text becomes "test" ---
$cols = "D", "F", "H", "I", "K", "L", "M", "N"
$rows = 2, 3, 4

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = "test"
    }
}

# Re-apply the style-1 look to each translated column (done per-column so no
# stray cells get created in the untouched E/G columns in between)
foreach ($col in $cols) {
    $styleSource.Copy()
    $ws.Range($col + "2:" + $col + "4").PasteSpecial(-4122)
}

# The "modified" date column (J) keeps its text but reverts to the default
# (unstyled) look instead of the style-1 look used by the translated cells
$ws.Range("J2:J4").ClearFormats()

$excel.CutCopyMode = $false

Write-Host "Applied he_IL translation updates"
